# Apply updated crypto price/volume data per commit
# "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.681.89'
$ws.Range('E2').Value = '  -3.44%  '
$ws.Range('D3').Value = '1.741.84'
$ws.Range('E3').Value = '  -5.78%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''236.10'
$ws.Range('E5').Value = '  -10.15%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '''0.4908'
$ws.Range('E7').Value = '  -8.30%  '
$ws.Range('D8').Value = '''41.37'
$ws.Range('E8').Value = '  -8.33%  '
$ws.Range('D9').Value = '''0.2553'
$ws.Range('E9').Value = '  -18.97%  '
$ws.Range('D10').Value = '''0.06030'
$ws.Range('E10').Value = '  -12.77%  '
$ws.Range('D11').Value = '1.743.58'
$ws.Range('E11').Value = '  -5.70%  '
$ws.Range('D12').Value = '''0.06829'
$ws.Range('E12').Value = '  -12.69%  '
$ws.Range('D13').Value = '''14.84'
$ws.Range('E13').Value = '  -21.06%  '
$ws.Range('D14').Value = '''4.434'
$ws.Range('E14').Value = '  -12.29%  '
$ws.Range('D15').Value = '''76.42'
$ws.Range('E15').Value = '  -14.70%  '
$ws.Range('D16').Value = '''0.5695'
$ws.Range('E16').Value = '  -25.88%  '
$ws.Range('D17').Value = '''1.001'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '25.723.48'
$ws.Range('E19').Value = '  -3.40%  '
$ws.Range('D20').Value = '''11.32'
$ws.Range('E20').Value = '  -19.72%  '
$ws.Range('D21').Value = '''0.000006556'
$ws.Range('E21').Value = '  -17.69%  '
$ws.Range('D22').Value = '1.965.74'
$ws.Range('E22').Value = '  -5.76%  '
$ws.Range('D23').Value = '''4.029'
$ws.Range('E23').Value = '  -13.19%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '''7.931'
$ws.Range('E24').Value = '  -15.12%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '''5.040'
$ws.Range('E25').Value = '  -16.23%  '
$ws.Range('D26').Value = '''136.83'
$ws.Range('E26').Value = '  -3.26%  '
$ws.Range('D27').Value = '''1.476'
$ws.Range('E27').Value = '  -12.58%  '
$ws.Range('D28').Value = '''1.814'
$ws.Range('E28').Value = '  -17.97%  '
$ws.Range('D29').Value = '''14.66'
$ws.Range('E29').Value = '  -13.91%  '
$ws.Range('E30').Value = '  -8.98%  '
$ws.Range('D31').Value = '''3.752'
$ws.Range('E31').Value = '  -12.80%  '
$ws.Range('D32').Value = '''0.07969'
$ws.Range('E32').Value = '  -9.25%  '
$ws.Range('D33').Value = '''3.388'
$ws.Range('E33').Value = '  -17.55%  '
$ws.Range('D34').Value = '''0.04385'
$ws.Range('E34').Value = '  -9.70%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').Value = '''2.611'
$ws.Range('E36').Value = '  -10.52%  '
$ws.Range('D37').Value = '''0.9812'
$ws.Range('E37').Value = '  -13.78%  '
$ws.Range('D38').Value = '''0.5981'
$ws.Range('E38').Value = '  -18.85%  '
$ws.Range('D39').Value = '''2.661'
$ws.Range('E39').Value = '  -14.50%  '
$ws.Range('D40').Value = '''1.924'
$ws.Range('E40').Value = '  -17.60%  '
$ws.Range('D41').Value = '''1.001'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = '''0.01511'
$ws.Range('E42').Value = '  -12.80%  '
$ws.Range('D43').Value = '''101.67'
$ws.Range('E43').Value = '  -6.28%  '
$ws.Range('D44').Value = '''0.7476'
$ws.Range('E44').Value = '  -17.47%  '
$ws.Range('D45').Value = '''5.155'
$ws.Range('E45').Value = '  -12.73%  '
$ws.Range('D46').Value = '''0.3737'
$ws.Range('E46').Value = '  -22.73%  '
$ws.Range('D47').Value = '''0.05228'
$ws.Range('E47').Value = '  -10.00%  '
$ws.Range('D48').Value = '''0.1066'
$ws.Range('E48').Value = '  -14.54%  '
$ws.Range('D49').Value = '''30.11'
$ws.Range('E49').Value = '  -14.14%  '
$ws.Range('D50').Value = '''52.08'
$ws.Range('E50').Value = '  -13.75%  '
$ws.Range('D51').Value = '''5.793'
$ws.Range('E51').Value = '  -24.64%  '
